# Apply the March 2025 attendance corrections for row 6 (R Rohit) and row 9 (Rajeev).
# Columns: M = Mar-16, N = Mar-17, O = Mar-18
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R Rohit's row (row 6): Mar-16 was blank -> "U", Mar-17 was "P**" -> "S", Mar-18 was "U" -> blank
$ws.Range("M6").Value = "U"
$ws.Range("N6").Value = "S"
$ws.Range("O6").Value = ""

# Rajeev's row (row 9): Mar-16 was "U" -> blank
$ws.Range("M9").Value = ""
